$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (shifts existing rows 59..148 down to 60..149)
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new data record
$ws.Cells.Item(59, 1).Value = 11
$ws.Cells.Item(59, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(59, 3).Value = "Bíobío"
$ws.Cells.Item(59, 4).Value = 45014
$ws.Cells.Item(59, 5).Value = 8
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100101
$ws.Cells.Item(59, 8).Value = "Berries"
$ws.Cells.Item(59, 9).Value = 100101001
$ws.Cells.Item(59, 10).Value = "Arándano (blue)"
$ws.Cells.Item(59, 11).Value = "Sin especificar"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 100
$ws.Cells.Item(59, 14).Value = 4500
$ws.Cells.Item(59, 15).Value = 5000
$ws.Cells.Item(59, 16).Value = 4750
$ws.Cells.Item(59, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(59, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(59, 19).Value = 2375
$ws.Cells.Item(59, 20).Value = 2
